# Split the "{m:userdoc 'zone1'}" field-like marker text into four
# separate runs: "{", "m", ":userdoc 'zone1'", "}".
#
# This mirrors the behaviour of the new TokenIteratorFieldRewriterSplit
# parser, which emits the pseudo-field token and its leading "{"/"m"
# marker and trailing "}" terminator as distinct runs instead of folding
# them into the two original runs.

$d = $word.ActiveDocument

# Locate the paragraph that holds the "{m:userdoc '...'}" marker text so
# this keeps working even if surrounding paragraphs shift around.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "{m:userdoc*") {
        $target = $para
        break
    }
}

$paraRange = $target.Range
$paraStart = $paraRange.Start

# The paragraph text is "{m:userdoc 'zone1'}" followed by the paragraph
# mark, laid out (initially) as two runs: "{m" and ":userdoc 'zone1'}".
#
#   offset 0  -> "{"
#   offset 1  -> "m"
#   offset 2  -> ":"
#   ...
#   offset 17 -> "'"
#   offset 18 -> "}"
#
# Toggling a character-formatting property on a sub-range and then right
# back to its original value forces the run to split at that boundary
# without altering the visible formatting.

# 1) Split "{m" -> "{" | "m"  (new run boundary after offset 1).
$splitA = $d.Range($paraStart, $paraStart + 1)
$splitA.Font.Bold = $true
$splitA.Font.Bold = $false

# 2) Split ":userdoc 'zone1'}" -> ":userdoc 'zone1'" | "}"
#    (new run boundary right before the closing "}", offset 18).
$splitB = $d.Range($paraStart + 2, $paraStart + 18)
$splitB.Font.Bold = $true
$splitB.Font.Bold = $false

Write-Output $target.Range.Text
